$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.056.19'
$ws.Range('E2').Value = '  +3.33%  '
$ws.Range('D3').Value = '2.448.01'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.18'
$ws.Range('E5').Value = '  +3.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.66'
$ws.Range('E6').Value = '  +2.79%  '
$ws.Range('E7').Value = '  +1.12%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  +4.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.01'
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0805'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.30'
$ws.Range('E13').Value = '  -2.58%  '
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('D15').Value = '2.830.48'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = '2.450.98'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.839'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '45.910.49'
$ws.Range('E18').Value = '  +3.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.54'
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.41'
$ws.Range('E21').Value = '  +2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.40'
$ws.Range('E22').Value = '  +3.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '247.61'
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.36'
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.03'
$ws.Range('E26').Value = '  +2.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.71'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.09'
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '49.38'
$ws.Range('E31').Value = '  +1.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.130'
$ws.Range('E32').Value = '  +6.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.92'
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.55'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '125.40'
$ws.Range('E40').Value = '  +0.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.27'
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.07'
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').Value = '1.958.57'
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('E46').Value = '  -0.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.82'
$ws.Range('E48').Value = '  +8.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.13'
$ws.Range('E49').Value = '  -6.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '77.73'
$ws.Range('E50').Value = '  +5.14%  '
$ws.Range('E51').Value = '  +5.66%  '
